$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.40"
$ws.Range("D3").Value = "'35.77"
$ws.Range("E3").Value = "'-0.42%"
$ws.Range("D4").Value = "'5.040"
$ws.Range("E4").Value = "'-0.33%"
$ws.Range("D5").Value = "'0.07969"
$ws.Range("E5").Value = "'-1.47%"
$ws.Range("D6").Value = "'1.870"
$ws.Range("E6").Value = "'-3.85%"
$ws.Range("D7").Value = "'7.774"
$ws.Range("E7").Value = "'-0.01%"
$ws.Range("D8").Value = "'0.9220"
$ws.Range("E8").Value = "'-0.88%"
$ws.Range("D9").Value = "'0.1278"
$ws.Range("E9").Value = "'-5.78%"
$ws.Range("D10").Value = "'0.1884"
$ws.Range("E10").Value = "'-1.13%"
$ws.Range("D11").Value = "'0.09077"
$ws.Range("E11").Value = "'-1.85%"
$ws.Range("E12").Value = "'-1.85%"
$ws.Range("D13").Value = "'0.09871"
$ws.Range("E13").Value = "'0.01%"
$ws.Range("D14").Value = "'0.001408"
$ws.Range("E14").Value = "'-0.20%"
$ws.Range("D15").Value = "'0.006201"
$ws.Range("E15").Value = "'7.04%"
$ws.Range("D16").Value = "'3.849"
$ws.Range("E16").Value = "'7.17%"
$ws.Range("D17").Value = "'4.128"
$ws.Range("E17").Value = "'-0.50%"
$ws.Range("D18").Value = "'3.280"
$ws.Range("E18").Value = "'10.23%"
$ws.Range("D19").Value = "'0.3409"
$ws.Range("E19").Value = "'-1.06%"
$ws.Range("D20").Value = "'0.1341"
$ws.Range("E20").Value = "'-0.30%"
$ws.Range("D21").Value = "'4.797"
$ws.Range("E21").Value = "'-1.83%"
$ws.Range("D22").Value = "'0.2504"
$ws.Range("E22").Value = "'-3.55%"
$ws.Range("D23").Value = "'0.04424"
$ws.Range("E23").Value = "'0.82%"
$ws.Range("E24").Value = "'1.14%"
$ws.Range("D25").Value = "'0.004856"
$ws.Range("E25").Value = "'1.72%"
$ws.Range("D27").Value = "'0.0001301"
$ws.Range("E27").Value = "'-21.19%"
$ws.Range("E28").Value = "'42.40%"
$ws.Range("D39").Value = "'0.01934"
$ws.Range("E39").Value = "'-2.37%"
$ws.Range("D40").Value = "'0.05161"
$ws.Range("E40").Value = "'1.60%"
$ws.Range("D41").Value = "'0.007542"
$ws.Range("E41").Value = "'-0.92%"
$ws.Range("D42").Value = "'0.01016"
$ws.Range("E42").Value = "'-9.32%"
$ws.Range("D43").Value = "'0.1345"
$ws.Range("E43").Value = "'-2.47%"
$ws.Range("D44").Value = "'0.002111"
$ws.Range("E44").Value = "'0.69%"
$ws.Range("D45").Value = "'0.009875"
$ws.Range("E45").Value = "'-8.73%"
$ws.Range("D46").Value = "'0.00006196"
$ws.Range("E46").Value = "'-2.75%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.24%"
$ws.Range("D48").Value = "'65.00"
$ws.Range("E48").Value = "'-0.34%"
$ws.Range("D49").Value = "'0.001253"
$ws.Range("E49").Value = "'5.45%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.24%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.24%"
